$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("B1").Value = "long edge"
$ws.Range("C1").Value = "length"
$ws.Range("D1").Value = "width"

# Move the active selection to D2
$ws.Range("D2").Select()
